# Append the new resale-number row (2025-02-03 11:31:58) to the
# CityResaleNum sheet, mirroring the existing data rows above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = 49

# Columns A and D look like dates/numbers to Excel's auto-detection
# ("2025-02-03" -> date serial, "05" -> number 5), but the source data
# keeps them as literal text, matching every other row in the sheet.
# Force text formatting before assigning, then drop back to the
# worksheet's default "Normal" style so no extra formatting sticks.
$ws.Cells.Item($r, 1).NumberFormat = "@"
$ws.Cells.Item($r, 1).Value = "2025-02-03"
$ws.Cells.Item($r, 1).Style = "Normal"

$ws.Cells.Item($r, 2).Value = "11:31:58"
$ws.Cells.Item($r, 3).Value = "Monday"

$ws.Cells.Item($r, 4).NumberFormat = "@"
$ws.Cells.Item($r, 4).Value = "05"
$ws.Cells.Item($r, 4).Style = "Normal"

$ws.Cells.Item($r, 5).Value = 125884
$ws.Cells.Item($r, 6).Value = 141882
$ws.Cells.Item($r, 7).Value = 166319
$ws.Cells.Item($r, 8).Value = 157817
$ws.Cells.Item($r, 9).Value = -1
$ws.Cells.Item($r, 10).Value = 142054
$ws.Cells.Item($r, 11).Value = -1
$ws.Cells.Item($r, 12).Value = -1
$ws.Cells.Item($r, 13).Value = 191120
$ws.Cells.Item($r, 14).Value = 115434
$ws.Cells.Item($r, 15).Value = 44744
$ws.Cells.Item($r, 16).Value = 28248
$ws.Cells.Item($r, 17).Value = 63150
$ws.Cells.Item($r, 18).Value = -1
$ws.Cells.Item($r, 19).Value = 39498
$ws.Cells.Item($r, 20).Value = -1
